$d = $word.ActiveDocument

# 1. Update the date line
$d.Content.Find.Execute("2023-09-25 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-26 Tuesday", 2) | Out-Null

$t = $d.Tables.Item(1)

# 2. Insert 6 new rows at the top, filled with the new practice problems.
#    Rows.Add(refRow) inserts immediately before refRow, so build them
#    in reverse order (last new row first) to end up in the right order.
$anchorRow = $t.Rows.Item(1)
$newRow = $t.Rows.Add($anchorRow)
$newRow.Cells.Item(1).Range.Text = "8+76=84"
$newRow.Cells.Item(2).Range.Text = "52-35=17"
$newRow.Cells.Item(3).Range.Text = "49-46=3"
$newRow.Cells.Item(4).Range.Text = "20-5=15"
$newRow.Cells.Item(5).Range.Text = "48-2=46"

$newRow = $t.Rows.Add($anchorRow)
$newRow.Cells.Item(1).Range.Text = "52+46=98"
$newRow.Cells.Item(2).Range.Text = "40+40=80"
$newRow.Cells.Item(3).Range.Text = "53-7=46"
$newRow.Cells.Item(4).Range.Text = "68-32=36"
$newRow.Cells.Item(5).Range.Text = "15+39=54"

$newRow = $t.Rows.Add($anchorRow)
$newRow.Cells.Item(1).Range.Text = "88-75=13"
$newRow.Cells.Item(2).Range.Text = "51+46=97"
$newRow.Cells.Item(3).Range.Text = "62-28=34"
$newRow.Cells.Item(4).Range.Text = "31-24=7"
$newRow.Cells.Item(5).Range.Text = "76-59=17"

$newRow = $t.Rows.Add($anchorRow)
$newRow.Cells.Item(1).Range.Text = "9+36=45"
$newRow.Cells.Item(2).Range.Text = "64-43=21"
$newRow.Cells.Item(3).Range.Text = "84+15=99"
$newRow.Cells.Item(4).Range.Text = "82-36=46"
$newRow.Cells.Item(5).Range.Text = "5+82=87"

$newRow = $t.Rows.Add($anchorRow)
$newRow.Cells.Item(1).Range.Text = "14+43=57"
$newRow.Cells.Item(2).Range.Text = "2+67=69"
$newRow.Cells.Item(3).Range.Text = "17+27=44"
$newRow.Cells.Item(4).Range.Text = "26-13=13"
$newRow.Cells.Item(5).Range.Text = "98-48=50"

$newRow = $t.Rows.Add($anchorRow)
$newRow.Cells.Item(1).Range.Text = "7+78=85"
$newRow.Cells.Item(2).Range.Text = "1+38=39"
$newRow.Cells.Item(3).Range.Text = "77-50=27"
$newRow.Cells.Item(4).Range.Text = "92-79=13"
$newRow.Cells.Item(5).Range.Text = "84-55=29"

# 3. Update the text of the (now shifted-down) pre-existing rows that are kept.
#    They are currently rows 7..20 (6 new rows + these 14).
$t.Rows.Item(7).Cells.Item(1).Range.Text = "98-71=27"
$t.Rows.Item(7).Cells.Item(2).Range.Text = "50-31=19"
$t.Rows.Item(7).Cells.Item(3).Range.Text = "16-2=14"
$t.Rows.Item(7).Cells.Item(4).Range.Text = "47+45=92"
$t.Rows.Item(7).Cells.Item(5).Range.Text = "53-3=50"

$t.Rows.Item(8).Cells.Item(1).Range.Text = "42+8=50"
$t.Rows.Item(8).Cells.Item(2).Range.Text = "36+63=99"
$t.Rows.Item(8).Cells.Item(3).Range.Text = "58-25=33"
$t.Rows.Item(8).Cells.Item(4).Range.Text = "6+22=28"
$t.Rows.Item(8).Cells.Item(5).Range.Text = "0+62=62"

$t.Rows.Item(9).Cells.Item(1).Range.Text = "32-25=7"
$t.Rows.Item(9).Cells.Item(2).Range.Text = "40-0=40"
$t.Rows.Item(9).Cells.Item(3).Range.Text = "88-49=39"
$t.Rows.Item(9).Cells.Item(4).Range.Text = "87-5=82"
$t.Rows.Item(9).Cells.Item(5).Range.Text = "25-16=9"

$t.Rows.Item(10).Cells.Item(1).Range.Text = "44+54=98"
$t.Rows.Item(10).Cells.Item(2).Range.Text = "77-1=76"
$t.Rows.Item(10).Cells.Item(3).Range.Text = "60+11=71"
$t.Rows.Item(10).Cells.Item(4).Range.Text = "84-80=4"
$t.Rows.Item(10).Cells.Item(5).Range.Text = "80-56=24"

$t.Rows.Item(11).Cells.Item(1).Range.Text = "25+23=48"
$t.Rows.Item(11).Cells.Item(2).Range.Text = "84-37=47"
$t.Rows.Item(11).Cells.Item(3).Range.Text = "30+66=96"
$t.Rows.Item(11).Cells.Item(4).Range.Text = "3+39=42"
$t.Rows.Item(11).Cells.Item(5).Range.Text = "73-12=61"

$t.Rows.Item(12).Cells.Item(1).Range.Text = "48+15=63"
$t.Rows.Item(12).Cells.Item(2).Range.Text = "78+0=78"
$t.Rows.Item(12).Cells.Item(3).Range.Text = "50+6=56"
$t.Rows.Item(12).Cells.Item(4).Range.Text = "43+43=86"
$t.Rows.Item(12).Cells.Item(5).Range.Text = "46+38=84"

$t.Rows.Item(13).Cells.Item(1).Range.Text = "16+30=46"
$t.Rows.Item(13).Cells.Item(2).Range.Text = "35-22=13"
$t.Rows.Item(13).Cells.Item(3).Range.Text = "7+2=9"
$t.Rows.Item(13).Cells.Item(4).Range.Text = "36+62=98"
$t.Rows.Item(13).Cells.Item(5).Range.Text = "15+59=74"

$t.Rows.Item(14).Cells.Item(1).Range.Text = "7+26=33"
$t.Rows.Item(14).Cells.Item(2).Range.Text = "84-28=56"
$t.Rows.Item(14).Cells.Item(3).Range.Text = "7+21=28"
$t.Rows.Item(14).Cells.Item(4).Range.Text = "16+15=31"
$t.Rows.Item(14).Cells.Item(5).Range.Text = "86-85=1"

$t.Rows.Item(15).Cells.Item(1).Range.Text = "64-28=36"
$t.Rows.Item(15).Cells.Item(2).Range.Text = "96-58=38"
$t.Rows.Item(15).Cells.Item(3).Range.Text = "16+32=48"
$t.Rows.Item(15).Cells.Item(4).Range.Text = "75-34=41"
$t.Rows.Item(15).Cells.Item(5).Range.Text = "48-40=8"

$t.Rows.Item(16).Cells.Item(1).Range.Text = "87-33=54"
$t.Rows.Item(16).Cells.Item(2).Range.Text = "15+41=56"
$t.Rows.Item(16).Cells.Item(3).Range.Text = "25+19=44"
$t.Rows.Item(16).Cells.Item(4).Range.Text = "32+25=57"
$t.Rows.Item(16).Cells.Item(5).Range.Text = "71-38=33"

$t.Rows.Item(17).Cells.Item(1).Range.Text = "95-49=46"
$t.Rows.Item(17).Cells.Item(2).Range.Text = "66-18=48"
$t.Rows.Item(17).Cells.Item(3).Range.Text = "43+48=91"
$t.Rows.Item(17).Cells.Item(4).Range.Text = "60+36=96"
$t.Rows.Item(17).Cells.Item(5).Range.Text = "1+63=64"

$t.Rows.Item(18).Cells.Item(1).Range.Text = "60-7=53"
$t.Rows.Item(18).Cells.Item(2).Range.Text = "78-6=72"
$t.Rows.Item(18).Cells.Item(3).Range.Text = "68-8=60"
$t.Rows.Item(18).Cells.Item(4).Range.Text = "84-53=31"
$t.Rows.Item(18).Cells.Item(5).Range.Text = "68+6=74"

$t.Rows.Item(19).Cells.Item(1).Range.Text = "46+4=50"
$t.Rows.Item(19).Cells.Item(2).Range.Text = "14+53=67"
$t.Rows.Item(19).Cells.Item(3).Range.Text = "97-13=84"
$t.Rows.Item(19).Cells.Item(4).Range.Text = "70-1=69"
$t.Rows.Item(19).Cells.Item(5).Range.Text = "44-13=31"

$t.Rows.Item(20).Cells.Item(1).Range.Text = "11+83=94"
$t.Rows.Item(20).Cells.Item(2).Range.Text = "61-26=35"
$t.Rows.Item(20).Cells.Item(3).Range.Text = "47+39=86"
$t.Rows.Item(20).Cells.Item(4).Range.Text = "8+56=64"
$t.Rows.Item(20).Cells.Item(5).Range.Text = "75-57=18"

# 4. Delete the trailing 6 rows that are no longer present in the target table.
for ($i = $t.Rows.Count; $i -gt (6 + 14); $i--) {
    $t.Rows.Item($i).Delete()
}

Write-Host "Final row count:" $t.Rows.Count
